$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.884.94'
$ws.Range("E2").Value = '  -0.74%  '

# Row 3
$ws.Range("D3").Value = '3.499.18'
$ws.Range("E3").Value = '  -1.95%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.87'
$ws.Range("E5").Value = '  -0.81%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '198.68'
$ws.Range("E6").Value = '  +6.43%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  +0.75%  '

# Row 8
$ws.Range("E8").Value = '  -0.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.210'
$ws.Range("E9").Value = '  -3.19%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.654'
$ws.Range("E10").Value = '  +1.12%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.28'
$ws.Range("E11").Value = '  +0.46%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000303'
$ws.Range("E12").Value = '  -2.58%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.58'
$ws.Range("E13").Value = '  +0.41%  '

# Row 14
$ws.Range("D14").Value = '4.056.12'
$ws.Range("E14").Value = '  -1.83%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '593.87'
$ws.Range("E15").Value = '  +2.27%  '

# Row 16
$ws.Range("D16").Value = '69.921.59'
$ws.Range("E16").Value = '  -0.80%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.03'
$ws.Range("E17").Value = '  -0.19%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.64'
$ws.Range("E18").Value = '  -1.24%  '

# Row 19
$ws.Range("D19").Value = '3.496.28'
$ws.Range("E19").Value = '  -2.90%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.121'
$ws.Range("E20").Value = '  +0.20%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.987'
$ws.Range("E21").Value = '  -1.01%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.06'
$ws.Range("E22").Value = '  +3.85%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '103.81'
$ws.Range("E23").Value = '  +9.70%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.59'
$ws.Range("E24").Value = '  -3.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.01'
$ws.Range("E25").Value = '  +2.83%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.09'
$ws.Range("E26").Value = '  +4.16%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.93'
$ws.Range("E27").Value = '  -0.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.75'
$ws.Range("E28").Value = '  +3.08%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.48'
$ws.Range("E29").Value = '  +3.41%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.51'
$ws.Range("E30").Value = '  +20.62%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.30'
$ws.Range("E31").Value = '  +3.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.73'
$ws.Range("E32").Value = '  +3.57%  '

# Row 33
$ws.Range("E33").Value = '  +0.72%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.56'
$ws.Range("E34").Value = '  -0.02%  '

# Row 35
$ws.Range("D35").Value = '3.722.98'
$ws.Range("E35").Value = '  +5.52%  '

# Row 36
$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").Value = '0.0₃0806'
$ws.Range("E36").Value = '  +2.20%  '

# Row 37
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.04%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '509.85'
$ws.Range("E38").Value = '  -4.37%  '

# Row 39
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.390'
$ws.Range("E39").Value = '  -3.82%  '

# Row 40
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.97'
$ws.Range("E40").Value = '  -7.84%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.63'
$ws.Range("E41").Value = '  -2.28%  '

# Row 42
$ws.Range("E42").Value = '  -0.79%  '

# Row 43
$ws.Range("E43").Value = '  +0.50%  '

# Row 44
$ws.Range("E44").Value = '  -1.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.84'
$ws.Range("E45").Value = '  -3.51%  '

# Row 46
$ws.Range("E46").Value = '  -1.35%  '

# Row 47
$ws.Range("E47").Value = '  -4.32%  '

# Row 48
$ws.Range("E48").Value = '  +0.27%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.76'
$ws.Range("E49").Value = '  -5.17%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.75'
$ws.Range("E50").Value = '  -2.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000241'
$ws.Range("E51").Value = '  -2.60%  '
